$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Y2").Value = 1545.169778625117
$ws.Range("Z2").Value = 7.302148999574287
$ws.Range("AA2").Value = 243.5881601706667
$ws.Range("AB2").Value = 1.151146666666667
$ws.Range("Y3").Value = 11959.50849732252
$ws.Range("Z3").Value = 17.32599193543705
$ws.Range("AA3").Value = 1617.946192544732
$ws.Range("AB3").Value = 2.343952737713009
$ws.Range("Y4").Value = 920.4965612866991
$ws.Range("Z4").Value = 10.06963511013957
$ws.Range("AA4").Value = 72.44918110161763
$ws.Range("AB4").Value = 0.7925470321170338
$ws.Range("Y5").Value = 53473.9996952143
$ws.Range("Z5").Value = 80.45496716249664
$ws.Range("AA5").Value = 11006.26021067463
$ws.Range("AB5").Value = 16.55960483372951
$ws.Range("Y6").Value = 4242.797033226667
$ws.Range("Z6").Value = 7.299983333333333
$ws.Range("AA6").Value = 544.4804411595852
$ws.Range("AB6").Value = 0.9368108148148148
$ws.Range("Y7").Value = 7206.456348488312
$ws.Range("Z7").Value = 10.01273583319321
$ws.Range("AA7").Value = 104.3261405805616
$ws.Range("AB7").Value = 0.1449519757861107
$ws.Range("Y8").Value = 16632.26992707227
$ws.Range("Z8").Value = 24.50405208569424
$ws.Range("AA8").Value = 1175.224044782467
$ws.Range("AB8").Value = 1.731438422535211
$ws.Range("Y9").Value = 1784.997211823123
$ws.Range("Z9").Value = 7.107122669220945
$ws.Range("AA9").Value = 161.3231664091966
$ws.Range("AB9").Value = 0.6423223103448277
$ws.Range("Y10").Value = 5335.778096408275
$ws.Range("Z10").Value = 8.266024137931034
$ws.Range("AA10").Value = 610.2831631103999
$ws.Range("AB10").Value = 0.9454319999999998
